$d = $word.ActiveDocument

$d.Content.Find.Execute("335×8=", $true, $false, $false, $false, $false, $true, 1, $false, "367×8=", 2) | Out-Null
$d.Content.Find.Execute("810×2=", $true, $false, $false, $false, $false, $true, 1, $false, "552×7=", 2) | Out-Null
$d.Content.Find.Execute("664×2=", $true, $false, $false, $false, $false, $true, 1, $false, "795×4=", 2) | Out-Null
$d.Content.Find.Execute("596×5=", $true, $false, $false, $false, $false, $true, 1, $false, "578×4=", 2) | Out-Null
$d.Content.Find.Execute("369×2=", $true, $false, $false, $false, $false, $true, 1, $false, "206×3=", 2) | Out-Null
$d.Content.Find.Execute("824×3=", $true, $false, $false, $false, $false, $true, 1, $false, "953×4=", 2) | Out-Null
$d.Content.Find.Execute("296×8=", $true, $false, $false, $false, $false, $true, 1, $false, "176×2=", 2) | Out-Null
$d.Content.Find.Execute("153×9=", $true, $false, $false, $false, $false, $true, 1, $false, "629×5=", 2) | Out-Null
$d.Content.Find.Execute("109×7=", $true, $false, $false, $false, $false, $true, 1, $false, "898×2=", 2) | Out-Null
$d.Content.Find.Execute("435×8=", $true, $false, $false, $false, $false, $true, 1, $false, "437×4=", 2) | Out-Null
$d.Content.Find.Execute("946×3=", $true, $false, $false, $false, $false, $true, 1, $false, "884×8=", 2) | Out-Null
$d.Content.Find.Execute("899×7=", $true, $false, $false, $false, $false, $true, 1, $false, "348×8=", 2) | Out-Null
$d.Content.Find.Execute("612×2=", $true, $false, $false, $false, $false, $true, 1, $false, "525×5=", 2) | Out-Null
$d.Content.Find.Execute("322×8=", $true, $false, $false, $false, $false, $true, 1, $false, "390×9=", 2) | Out-Null
$d.Content.Find.Execute("810×8=", $true, $false, $false, $false, $false, $true, 1, $false, "994×5=", 2) | Out-Null
$d.Content.Find.Execute("965×5=", $true, $false, $false, $false, $false, $true, 1, $false, "539×6=", 2) | Out-Null
$d.Content.Find.Execute("250×7=", $true, $false, $false, $false, $false, $true, 1, $false, "643×5=", 2) | Out-Null
$d.Content.Find.Execute("627×9=", $true, $false, $false, $false, $false, $true, 1, $false, "235×4=", 2) | Out-Null
$d.Content.Find.Execute("232×9=", $true, $false, $false, $false, $false, $true, 1, $false, "448×3=", 2) | Out-Null
$d.Content.Find.Execute("225×8=", $true, $false, $false, $false, $false, $true, 1, $false, "515×4=", 2) | Out-Null
$d.Content.Find.Execute("279×4=", $true, $false, $false, $false, $false, $true, 1, $false, "338×6=", 2) | Out-Null
$d.Content.Find.Execute("704×6=", $true, $false, $false, $false, $false, $true, 1, $false, "495×8=", 2) | Out-Null
$d.Content.Find.Execute("267×7=", $true, $false, $false, $false, $false, $true, 1, $false, "879×9=", 2) | Out-Null
$d.Content.Find.Execute("866×7=", $true, $false, $false, $false, $false, $true, 1, $false, "229×4=", 2) | Out-Null
$d.Content.Find.Execute("742×7=", $true, $false, $false, $false, $false, $true, 1, $false, "918×3=", 2) | Out-Null
